$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the promotional text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$oldText = $cellA1.Value2
$search = "✅ 1000 Bs = 2.32 = 8863.11 pesos`n✅ 8863.11 pesos = 2.31 = 959.36 Bs"
$replacement = "✅ 1000 Bs = 2.35 = 8999.76 pesos`n✅ 8999.76 pesos = 2.37 = 963.93 Bs"
if (-not $oldText.Contains($search)) {
    throw "Expected rate lines not found in A1 text; aborting to avoid silent no-op."
}
$newText = $oldText.Replace($search, $replacement)
$cellA1.Value2 = $newText

# --- Sheet "tasas": update rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value2 = 424.7
$ws2.Range("O10").Value2 = 3822.2
$ws2.Range("N12").Value2 = 3800
$ws2.Range("O12").Value2 = 407.002
